# Insert a new row at position 55 (pushes existing rows 55..186 down to 56..187)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with its data
$ws.Cells.Item(55, 1).Value = 3
$ws.Cells.Item(55, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(55, 3).Value = "Coquimbo"
$ws.Cells.Item(55, 4).Value = 44526
$ws.Cells.Item(55, 5).Value = 5
$ws.Cells.Item(55, 6).Value = 100112001
$ws.Cells.Item(55, 7).Value = "Berenjena"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 40
$ws.Cells.Item(55, 11).Value = 7000
$ws.Cells.Item(55, 12).Value = 7000
$ws.Cells.Item(55, 13).Value = 7000
$ws.Cells.Item(55, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 117
$ws.Cells.Item(55, 17).Value = 60
$ws.Cells.Item(55, 18).Value = "Hortaliza"
